$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("10:10").Copy()
$ws.Rows("11:11").Insert()
Write-Output "done"
